# Update "想去人数" (column F) values for a handful of rows on both the
# "展览" sheet and the "全部类型" sheet (which mirrors the same data).

$wb = $excel.ActiveWorkbook

# Row number -> new value for column F
$updates = @{
    4  = 353
    8  = 120
    11 = 62
    13 = 1103
    14 = 1454
    18 = 94
    22 = 252
    23 = 276
    25 = 1672
    29 = 627
    30 = 300
    31 = 3932
    32 = 9
    34 = 227
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
